# Re-applies the raw-data corrections from the commit "updates to raw data and
# Stata codes" to analytics/modified_data/tbf_market_garden_sales_2023_clean.xlsx.
#
# The edits are plain cell-value overwrites on Sheet1: a handful of
# `sale_item` labels were corrected/renamed (columns G "sale_item" and H
# "sale_item_code" always carry the same text) and many `sale_amnt` (I) /
# `sale_value_usd` (L) figures were corrected. Setting `.Value` lets Excel's
# COM layer manage the shared-string table itself, so we only need the
# final resolved text/number for each touched cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = "G11"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H11"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I11"; Value = 2 },
    @{ Cell = "G12"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H12"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I12"; Value = 2 },
    @{ Cell = "G13"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H13"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I13"; Value = 2 },
    @{ Cell = "G14"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H14"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I14"; Value = 2 },
    @{ Cell = "G15"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H15"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I15"; Value = 4 },
    @{ Cell = "G16"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H16"; Value = 'pepper- califonia wonder' },
    @{ Cell = "L16"; Value = 2 },
    @{ Cell = "G17"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H17"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I17"; Value = 1 },
    @{ Cell = "L17"; Value = 1 },
    @{ Cell = "G18"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H18"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I18"; Value = 5 },
    @{ Cell = "G19"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H19"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I19"; Value = 5 },
    @{ Cell = "G20"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H20"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I20"; Value = 6 },
    @{ Cell = "G21"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H21"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I21"; Value = 10 },
    @{ Cell = "L21"; Value = 2 },
    @{ Cell = "G22"; Value = 'pepper- sweet banana' },
    @{ Cell = "H22"; Value = 'pepper- sweet banana' },
    @{ Cell = "L22"; Value = 1 },
    @{ Cell = "G23"; Value = 'pepper- sweet banana' },
    @{ Cell = "H23"; Value = 'pepper- sweet banana' },
    @{ Cell = "I23"; Value = 8 },
    @{ Cell = "G37"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H37"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I37"; Value = 2 },
    @{ Cell = "G38"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H38"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I38"; Value = 2 },
    @{ Cell = "G39"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H39"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I39"; Value = 4 },
    @{ Cell = "L39"; Value = 2 },
    @{ Cell = "G40"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H40"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I40"; Value = 6 },
    @{ Cell = "G41"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H41"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I41"; Value = 5 },
    @{ Cell = "G42"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H42"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I42"; Value = 5 },
    @{ Cell = "L42"; Value = 1 },
    @{ Cell = "G43"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H43"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I43"; Value = 5 },
    @{ Cell = "G56"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H56"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "G57"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H57"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "G66"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H66"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I66"; Value = 2 },
    @{ Cell = "G67"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H67"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I67"; Value = 4 },
    @{ Cell = "G68"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H68"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I68"; Value = 7 },
    @{ Cell = "L68"; Value = 3 },
    @{ Cell = "G69"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H69"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I69"; Value = 4 },
    @{ Cell = "L69"; Value = 1 },
    @{ Cell = "G70"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H70"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I70"; Value = 8 },
    @{ Cell = "L70"; Value = 2 },
    @{ Cell = "G71"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H71"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I71"; Value = 20 },
    @{ Cell = "L71"; Value = 5 },
    @{ Cell = "G72"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H72"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I72"; Value = 55 },
    @{ Cell = "L72"; Value = 10 },
    @{ Cell = "G73"; Value = 'pepper- sweet banana' },
    @{ Cell = "H73"; Value = 'pepper- sweet banana' },
    @{ Cell = "L73"; Value = 1 },
    @{ Cell = "G74"; Value = 'pepper- sweet banana' },
    @{ Cell = "H74"; Value = 'pepper- sweet banana' },
    @{ Cell = "I74"; Value = 10 },
    @{ Cell = "L74"; Value = 2 },
    @{ Cell = "G93"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H93"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I93"; Value = 1 },
    @{ Cell = "G94"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H94"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I94"; Value = 2 },
    @{ Cell = "L94"; Value = 1 },
    @{ Cell = "G95"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H95"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I95"; Value = 2 },
    @{ Cell = "L95"; Value = 1 },
    @{ Cell = "G96"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H96"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I96"; Value = 2 },
    @{ Cell = "L96"; Value = 1 },
    @{ Cell = "G97"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H97"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I97"; Value = 2 },
    @{ Cell = "L97"; Value = 0.5 },
    @{ Cell = "G98"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H98"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "L98"; Value = 3 },
    @{ Cell = "G99"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H99"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I99"; Value = 8 },
    @{ Cell = "L99"; Value = 2 },
    @{ Cell = "G100"; Value = 'pepper- sweet banana' },
    @{ Cell = "H100"; Value = 'pepper- sweet banana' },
    @{ Cell = "I100"; Value = 3 },
    @{ Cell = "L100"; Value = 0.5 },
    @{ Cell = "G101"; Value = 'pepper- sweet banana' },
    @{ Cell = "H101"; Value = 'pepper- sweet banana' },
    @{ Cell = "I101"; Value = 3 },
    @{ Cell = "G102"; Value = 'pepper- sweet banana' },
    @{ Cell = "H102"; Value = 'pepper- sweet banana' },
    @{ Cell = "I102"; Value = 5 },
    @{ Cell = "G121"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H121"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I121"; Value = 1 },
    @{ Cell = "L121"; Value = 0.5 },
    @{ Cell = "G122"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H122"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I122"; Value = 1 },
    @{ Cell = "G123"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H123"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I123"; Value = 4 },
    @{ Cell = "G124"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H124"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I124"; Value = 4 },
    @{ Cell = "G125"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H125"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I125"; Value = 5 },
    @{ Cell = "L125"; Value = 1 },
    @{ Cell = "G126"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H126"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I126"; Value = 6 },
    @{ Cell = "I127"; Value = 2 },
    @{ Cell = "L127"; Value = 0.5 },
    @{ Cell = "I128"; Value = 5 },
    @{ Cell = "G129"; Value = 'pepper- sweet banana' },
    @{ Cell = "H129"; Value = 'pepper- sweet banana' },
    @{ Cell = "I129"; Value = 5 },
    @{ Cell = "L129"; Value = 1 },
    @{ Cell = "G130"; Value = 'pepper- sweet banana' },
    @{ Cell = "H130"; Value = 'pepper- sweet banana' },
    @{ Cell = "I130"; Value = 8 },
    @{ Cell = "G131"; Value = 'pumpkin- connecticut field' },
    @{ Cell = "H131"; Value = 'pumpkin- connecticut field' },
    @{ Cell = "G149"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H149"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "G157"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H157"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I157"; Value = 2 },
    @{ Cell = "G158"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H158"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I158"; Value = 2 },
    @{ Cell = "L158"; Value = 1 },
    @{ Cell = "G159"; Value = 'pepper- califonia wonder' },
    @{ Cell = "H159"; Value = 'pepper- califonia wonder' },
    @{ Cell = "I159"; Value = 4 },
    @{ Cell = "G160"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H160"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I160"; Value = 5 },
    @{ Cell = "G161"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "H161"; Value = 'pepper- jalapeno, early' },
    @{ Cell = "I161"; Value = 10 },
    @{ Cell = "L161"; Value = 2 },
    @{ Cell = "G162"; Value = 'pepper- sweet banana' },
    @{ Cell = "H162"; Value = 'pepper- sweet banana' },
    @{ Cell = "I162"; Value = 11 },
    @{ Cell = "G163"; Value = 'pumpkin- connecticut field' },
    @{ Cell = "H163"; Value = 'pumpkin- connecticut field' },
    @{ Cell = "G177"; Value = 'corn- ornamental, glass gem' },
    @{ Cell = "H177"; Value = 'corn- ornamental, glass gem' }
)

foreach ($edit in $edits) {
    $ws.Range($edit.Cell).Value = $edit.Value
}
